$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 244; this shifts the existing rows 244-299
# down to 245-300, matching the diff (row 299's old content ends up as
# the new row 300).
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new record's data.
$ws.Range("A244").Value = 5
$ws.Range("B244").Value = "Macroferia Regional de Talca"
$ws.Range("C244").Value = "Maule"
$ws.Range("D244").Value = 45275
$ws.Range("E244").Value = 7
$ws.Range("F244").Value = 100112031
$ws.Range("G244").Value = "Poroto verde"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 100
$ws.Range("K244").Value = 25000
$ws.Range("L244").Value = 28000
$ws.Range("M244").Value = 26500
$ws.Range("N244").Value = "$/saco 25 kilos"
$ws.Range("O244").Value = "Región del Maule"
$ws.Range("P244").Value = 1060
$ws.Range("Q244").Value = 25
$ws.Range("R244").Value = "Hortaliza"
